$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-CellText($row, $col, $old, $new) {
    $c = $t.Cell($row, $col)
    $r = $d.Range($c.Range.Start, $c.Range.End)
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

Replace-CellText 1 1 "71÷9=" "79÷5="
Replace-CellText 1 2 "41÷4=" "50÷5="
Replace-CellText 1 3 "87÷4=" "30÷3="
Replace-CellText 1 4 "99÷4=" "17÷7="
Replace-CellText 1 5 "87÷6=" "90÷9="

Replace-CellText 5 1 "95÷9=" "68÷7="
Replace-CellText 5 2 "96÷9=" "81÷9="
Replace-CellText 5 3 "99÷4=" "57÷2="
Replace-CellText 5 4 "94÷9=" "10÷5="
Replace-CellText 5 5 "67÷5=" "89÷5="

Replace-CellText 9 1 "42÷6=" "31÷2="
Replace-CellText 9 2 "96÷8=" "57÷6="
Replace-CellText 9 3 "26÷6=" "71÷8="
Replace-CellText 9 4 "33÷3=" "58÷2="
Replace-CellText 9 5 "61÷7=" "48÷8="

Replace-CellText 13 1 "72÷3=" "55÷2="
Replace-CellText 13 2 "92÷4=" "79÷6="
Replace-CellText 13 3 "41÷5=" "49÷5="
Replace-CellText 13 4 "88÷8=" "96÷9="
Replace-CellText 13 5 "69÷8=" "79÷9="

Replace-CellText 17 1 "60÷3=" "91÷3="
Replace-CellText 17 2 "27÷2=" "98÷5="
Replace-CellText 17 3 "48÷8=" "76÷5="
Replace-CellText 17 4 "67÷2=" "59÷3="
Replace-CellText 17 5 "29÷5=" "99÷6="
